$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "the opportunity to apply technical skills in real-world situations, or the potential for career growth and development.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "the opportunity to apply technical skills to real-world problems, or the potential for career growth and development.",
    2)

$d.Content.Find.Execute(
    "and continuous monitoring. We also prioritize safety training and education to ensure that everyone is aware of the potential hazards and knows how to respond in case of an emergency.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and regular training. We also conduct regular safety audits and inspections to identify any potential hazards and take corrective actions as needed. Additionally, we encourage a safety culture where everyone is responsible for maintaining a safe work environment and reporting any safety concerns or incidents promptly.",
    2)

$d.Content.Find.Execute(
    "to ensure pipeline integrity and flow efficiency, it's essential to conduct regular inspections, maintenance, and cleaning of the pipelines using appropriate pigging techniques. additionally, monitoring the pipeline's performance and analyzing data can help identify potential issues and optimize the pigging process.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "to ensure pipeline integrity and flow efficiency, it's essential to conduct regular inspections, maintenance, and cleaning of pipelines using appropriate pigging techniques. additionally, monitoring pipeline performance and analyzing data can help identify any potential issues and optimize pipeline operations.",
    2)

$d.Content.Find.Execute(
    "managing pressure control devices and conducting choke changes requires a deep understanding of the production process, equipment, and safety protocols. to optimize production rates and ensure safety, it's crucial to monitor the pressure and flow rates, identify potential issues, and adjust the pressure control devices and chokes accordingly. additionally, regular maintenance and inspection of the equipment can help prevent failures and ensure safe and efficient operations.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "managing pressure control devices and conducting choke changes requires a deep understanding of wellbore and reservoir dynamics, as well as the ability to operate and maintain various types of pressure control equipment. to optimize production rates and ensure safety, it's crucial to monitor well performance, analyze data, and make informed decisions based on the results. additionally, following proper procedures and safety protocols is essential to minimize risks and prevent incidents.",
    2)
